$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update existing summary figures (VALOR MORA / counts)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 1020000    # VALOR MORA total
$ws.Range("C13").Value2 = 3          # Cant. Trabajadores
$ws.Range("F13").Value2 = 5          # Cant. Periodos

# ---------------------------------------------------------------------------
# 2) Bump the "Valor Mora" amounts on the existing detail rows (16-18) and
#    on the row that is about to become a normal (non-last) row (19)
# ---------------------------------------------------------------------------
$ws.Range("G16").Value2 = 3500000
$ws.Range("G17").Value2 = 3500000
$ws.Range("G18").Value2 = 3500000
$ws.Range("G19").Value2 = 3500000

# ---------------------------------------------------------------------------
# 3) Make room for 3 new detail rows right after the table (which currently
#    ends at row 19). Inserting whole rows 20-22 pushes everything below
#    (the signature footer, rows 24-25) down to rows 27-28, correctly
#    relocating both the cell content and the merged-cell ranges that cover
#    it (B24:C24 -> B27:C27, etc.)
# ---------------------------------------------------------------------------
$ws.Rows("20:22").Insert()

# ---------------------------------------------------------------------------
# 4) Turn row 19 into a normal interior table row and grow the table with
#    2 more interior rows (20-21) plus a new closing row (22). Preserve the
#    special "closing" (thick-bottom-border) formatting that used to live
#    on row 19 by copying it onto the new last row (22) first, then
#    restyle row 19 (and the new rows 20-21) from row 18's normal style.
# ---------------------------------------------------------------------------
$ws.Range("B19:J19").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B20:J20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Fill in row 19 (keeps its original worker/period - only "Valor Mora"
#    changed, already applied above) and the 3 new worker/period detail
#    rows (20, 21, 22).
# ---------------------------------------------------------------------------

# Row 20: new worker MARILUZ HERNANDEZ ARIAS, period 2508
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1047394818"
$ws.Range("D20").Value2 = "MARILUZ HERNANDEZ ARIAS"
$ws.Range("E20").Value2 = "2508"
$ws.Range("F20").Value2 = 200000
$ws.Range("G20").Value2 = 5000000

# Row 21: existing worker EDWIN GUILLERMO HERNANDEZ ARIAS, period 2508
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1143326450"
$ws.Range("D21").Value2 = "EDWIN GUILLERMO HERNANDEZ ARIAS"
$ws.Range("E21").Value2 = "2508"
$ws.Range("F21").Value2 = 140000
$ws.Range("G21").Value2 = 3500000

# Row 22 (closing row): new worker STEEL POLO VARGAS, period 2508
$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1007446090"
$ws.Range("D22").Value2 = "STEEL POLO VARGAS"
$ws.Range("E22").Value2 = "2508"
$ws.Range("F22").Value2 = 120000
$ws.Range("G22").Value2 = 3000000

Write-Host "Edit complete"
